# Series para el dia 30 de Marzo de 2020
# Adds a new data row (row 28) to the "casos_chile_regiones" sheet with the
# case counts for the new day, and updates the view state (selection /
# scroll position) to match where the user ended up after typing the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 28

# fecha (date serial) - keep the existing DD/MM/YY date format used by
# column A (same style as the row above it).
$ws.Cells.Item($newRow, 1).Value = 43919
$ws.Cells.Item($newRow, 1).NumberFormat = "DD/MM/YY"

# dia, then the 16 region columns, then the total.
$values = @(27, 6, 8, 35, 2, 27, 108, 1295, 21, 54, 229, 201, 247, 47, 128, 2, 39, 2449)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($newRow, 2 + $i).Value = $values[$i]
}

# Leave the view where the user's cursor landed after entering the row.
[void]$ws.Range("N39").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 7
